$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so the refreshed model values can be written.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer note (A59).
$ws.Cells.Item(59, 1).Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."
$ws.Rows.Item(59).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns for each holding row with the latest model data.

$ws.Cells.Item(2, 4).Value = 0.02371495518332394
$ws.Cells.Item(2, 5).Value = -0.0002876042565430614
$ws.Cells.Item(3, 4).Value = 0.01769633165041824
$ws.Cells.Item(3, 5).Value = 0.001585728444003953
$ws.Cells.Item(4, 4).Value = 0.0182914709176837
$ws.Cells.Item(4, 5).Value = -0.003323958066990595
$ws.Cells.Item(5, 4).Value = 0.02081136150871507
$ws.Cells.Item(5, 5).Value = -0.008708272859216382
$ws.Cells.Item(6, 4).Value = 0.01977873446967128
$ws.Cells.Item(6, 5).Value = -0.005675100496571384
$ws.Cells.Item(7, 4).Value = 0.0274179140427144
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(8, 4).Value = 0.01903159499662092
$ws.Cells.Item(8, 5).Value = 0.005877413937867448
$ws.Cells.Item(9, 4).Value = 0.0198722730578466
$ws.Cells.Item(9, 5).Value = -0.01788656154389268
$ws.Cells.Item(10, 4).Value = 0.01918710289946238
$ws.Cells.Item(10, 5).Value = -0.006459475929311287
$ws.Cells.Item(11, 4).Value = 0.01981264220788484
$ws.Cells.Item(11, 5).Value = -0.01047506639126594
$ws.Cells.Item(12, 4).Value = 0.01950747256396285
$ws.Cells.Item(12, 5).Value = -0.01159793814432986
$ws.Cells.Item(13, 4).Value = 0.01976431393732759
$ws.Cells.Item(13, 5).Value = 0.02212537713711016
$ws.Cells.Item(14, 4).Value = 0.01883009728792658
$ws.Cells.Item(14, 5).Value = 0.01434366850188362
$ws.Cells.Item(15, 4).Value = 0.01768892651218769
$ws.Cells.Item(15, 5).Value = 0.002423654871546388
$ws.Cells.Item(16, 4).Value = 0.01771815732099248
$ws.Cells.Item(16, 5).Value = 0.02201887331998864
$ws.Cells.Item(17, 4).Value = 0.01572325205609509
$ws.Cells.Item(17, 5).Value = -0.008179959100204526
$ws.Cells.Item(18, 4).Value = 0.01566030838113545
$ws.Cells.Item(18, 5).Value = 0.02346880366342297
$ws.Cells.Item(19, 4).Value = 0.0165335300761638
$ws.Cells.Item(19, 5).Value = 0.00907558667185282
$ws.Cells.Item(20, 4).Value = 0.01965479584033899
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 4).Value = 0.01933793387289509
$ws.Cells.Item(21, 5).Value = 0.01795755487030659
$ws.Cells.Item(22, 4).Value = 0.02108983368059534
$ws.Cells.Item(22, 5).Value = -0.01480263157894735
$ws.Cells.Item(23, 4).Value = 0.01866971758361765
$ws.Cells.Item(23, 5).Value = -0.01221230624706426
$ws.Cells.Item(24, 4).Value = 0.02060966892796204
$ws.Cells.Item(24, 5).Value = 0.001701966717095216
$ws.Cells.Item(25, 4).Value = 0.02077492043373843
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 0.01923387219355005
$ws.Cells.Item(26, 5).Value = 0.01756838905775071
$ws.Cells.Item(27, 4).Value = 0.01965830353739556
$ws.Cells.Item(27, 5).Value = -0.0246238030095759
$ws.Cells.Item(28, 4).Value = 0.0273758216780355
$ws.Cells.Item(28, 5).Value = 0.009866102889358741
$ws.Cells.Item(29, 4).Value = 0.01884217935556589
$ws.Cells.Item(29, 5).Value = 0.002171889543903127
$ws.Cells.Item(30, 4).Value = 0.01281868402119273
$ws.Cells.Item(30, 5).Value = 0.003952569169960229
$ws.Cells.Item(31, 4).Value = 0.009439602523359313
$ws.Cells.Item(31, 5).Value = -0.009393063583814976
$ws.Cells.Item(32, 4).Value = 0.01671924314810355
$ws.Cells.Item(32, 5).Value = -0.005186721991701337
$ws.Cells.Item(33, 4).Value = 0.01999426296659191
$ws.Cells.Item(33, 5).Value = -0.01356698699830416
$ws.Cells.Item(34, 4).Value = 0.01867225092038073
$ws.Cells.Item(34, 5).Value = -0.000480076812289898
$ws.Cells.Item(35, 4).Value = 0.01899184109664641
$ws.Cells.Item(35, 5).Value = 0.02809415337889121
$ws.Cells.Item(36, 4).Value = 0.01785729597090326
$ws.Cells.Item(36, 5).Value = -0.006089309878213878
$ws.Cells.Item(37, 4).Value = 0.01959867268743379
$ws.Cells.Item(37, 5).Value = 0.002028397565923212
$ws.Cells.Item(38, 4).Value = 0.01949597511249964
$ws.Cells.Item(38, 5).Value = 0.002748763056624304
$ws.Cells.Item(39, 4).Value = 0.02513381864270832
$ws.Cells.Item(39, 5).Value = 0.005954596203944895
$ws.Cells.Item(40, 4).Value = 0.01655788908350112
$ws.Cells.Item(40, 5).Value = 0.004472271914132353
$ws.Cells.Item(41, 4).Value = 0.02169783450373492
$ws.Cells.Item(41, 5).Value = 0.007472337979594679
$ws.Cells.Item(42, 4).Value = 0.01912376948038535
$ws.Cells.Item(42, 5).Value = 0.006776379477250627
$ws.Cells.Item(43, 4).Value = 0.02013125802385702
$ws.Cells.Item(43, 5).Value = 0.008247422680412342
$ws.Cells.Item(44, 4).Value = 0.01820611695597373
$ws.Cells.Item(44, 5).Value = -0.009354997538158383
$ws.Cells.Item(45, 4).Value = 0.02008156564888888
$ws.Cells.Item(45, 5).Value = 0.01431344007763213
$ws.Cells.Item(46, 4).Value = 0.0192537491435373
$ws.Cells.Item(46, 5).Value = 0.008684034736138768
$ws.Cells.Item(47, 4).Value = 0.01789568576646688
$ws.Cells.Item(47, 5).Value = 0.01715069746169684
$ws.Cells.Item(48, 4).Value = 0.01613540646024259
$ws.Cells.Item(48, 5).Value = 0.007246376811594235
$ws.Cells.Item(49, 4).Value = 0.0180447628913713
$ws.Cells.Item(49, 5).Value = 0.01117734724292108
$ws.Cells.Item(50, 4).Value = 0.01668260720106821
$ws.Cells.Item(50, 5).Value = 0.04709840201850279
$ws.Cells.Item(51, 4).Value = 0.01669059695547485
$ws.Cells.Item(51, 5).Value = -0.01340354236476782
$ws.Cells.Item(52, 4).Value = 0.01823612725301331
$ws.Cells.Item(52, 5).Value = 0.01282325283180175
$ws.Cells.Item(53, 4).Value = 0.01545764144008892
$ws.Cells.Item(53, 5).Value = 0.01157308186883865
$ws.Cells.Item(54, 4).Value = 0.007629241098049486
$ws.Cells.Item(54, 5).Value = 0.01417624521072791
$ws.Cells.Item(55, 4).Value = 0.007166614830699053
$ws.Cells.Item(55, 5).Value = 0.003480530780944235
$ws.Cells.Item(56, 4).Value = 0.9999999999999999
$ws.Cells.Item(56, 5).Value = 0.003214414608233351

# Restore sheet protection.
$ws.Protect()
